# Bind the PathFile column (H) into the Item table and update the
# active selection to reflect it, per commit:
# "[Table] ItemTable, BasePathTable 수정, 아이템 에셋 - 테이블 바인딩"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header rows for the new "PathFile" column (H)
$ws.Range("H3").Value = "PathFile"
$ws.Range("H4").Value = "int32"

# Per-row PathFile ids bound to each item asset
$ws.Range("H5").Value  = 5006
$ws.Range("H6").Value  = 5004
$ws.Range("H7").Value  = 5005
$ws.Range("H8").Value  = 5008
$ws.Range("H9").Value  = 5009
$ws.Range("H10").Value = 5007
$ws.Range("H11").Value = 5001
$ws.Range("H12").Value = 5002
$ws.Range("H13").Value = 5010
$ws.Range("H14").Value = 5011
$ws.Range("H15").Value = 5012
$ws.Range("H16").Value = 5003

# Reflect the newly-filled column as the active selection
$ws.Range("H3:H16").Select()
